# Generate Report for Handoff
# b.md has been handed off again (new handback source produced a new
# handoff target b.63290e5768f688058c7b37413b0a5c26c308f864.*), so the
# "Overview" summary row and the per-locale detail rows for b.md move
# from "Handed back: in sync with en-US" to "Ready for handoff" with a
# fresh handoff file name + timestamp.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------
# Overview sheet: row 3 is the b.md summary row
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady
$wsOverview.Range("D3").Value = "2016-03-19 16:41:56"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md detail row
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 16:41:49"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------
# de-de sheet: row 3 is the b.md detail row
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 16:41:56"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
